$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Update existing ConvTranspose block (A/B) ---
$ws.Range("B2").Value = 188

# --- Update existing Conv block (E/F) ---
$ws.Range("F2").Value = 47
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 2

# --- Add new Maxpool block (H/I) ---
$ws.Range("H1").Value = "Maxpool"

$ws.Range("H2").Value = "Input"
$ws.Range("I2").Value = 94
$ws.Range("H3").Value = "Maxpool"
$ws.Range("I3").Value = 2
$ws.Range("H4").Value = "Dimension"
$ws.Range("I4").Formula = "=I2/I3"

$ws.Range("H6").Value = "Input"
$ws.Range("I6").Formula = "=I4"
$ws.Range("H7").Value = "Maxpool"
$ws.Range("I7").Value = 2
$ws.Range("H8").Value = "Dimension"
$ws.Range("I8").Formula = "=I6/I7"

$ws.Range("H10").Value = "Input"
$ws.Range("I10").Formula = "=I8"
$ws.Range("H11").Value = "Maxpool"
$ws.Range("I11").Value = 2
$ws.Range("H12").Value = "Dimension"
$ws.Range("I12").Formula = "=I10/I11"

$ws.Range("H14").Value = "Input"
$ws.Range("I14").Formula = "=I12"
$ws.Range("H15").Value = "Maxpool"
$ws.Range("I15").Value = 2
$ws.Range("H16").Value = "Dimension"
$ws.Range("I16").Formula = "=I14/I15"

# --- Selection matches target workbook state ---
$ws.Range("C15").Select()

$wb.Save()
